$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "target" column (G) changes its value from "helium" to "h" for every
# data row (rows 2-10). Shared strings are recomputed automatically by the
# engine: once no cell references "helium" anymore it drops out of the
# sharedStrings table and "h" is appended, matching the diff.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 7).Value = "h"
}

# The header row (row 1, columns A:K) becomes bold and horizontally
# centered, which introduces a new font + cell style in styles.xml.
$headerRange = $ws.Range("A1:K1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter

# The active selection moved from K11 to F15.
$ws.Range("F15").Select()
